# Apply repulled "dSF" (column F) values for specific rows in the
# save-data log sheet. These are data corrections from a repull of the
# underlying data source (see commit message: "repull data, push all
# data, mean calculation").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -4
    3  = -2
    4  = -1
    6  = -4
    8  = -2
    11 = 2
    12 = -3
    13 = 3
    14 = -3
    15 = -2
    16 = 3
    17 = -2
    18 = -2
    19 = -4
    23 = -1
    27 = -7
    31 = -4
    32 = -8
    38 = -9
    39 = -6
    42 = -10
    43 = -4
    45 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
